# Generate Report for Handoff
#
# The 9000b78f-fc35-4fc7-ad39-fc29d777b894.md file has finished translation
# and is now ready to be handed off. Update its status/priority/timestamps
# on the per-language sheets (zh-cn, de-de) and roll the summary up onto the
# Overview sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newHandoffDateTime = "2016-08-21 18:11:45"
$newGenerateDate     = "2016-08-21 18:11:49"

# --- zh-cn sheet: row 3 is the 9000b78f...md file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = $newHandoffDateTime

# --- de-de sheet: row 3 is the 9000b78f...md file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = $newGenerateDate

# --- Overview sheet: row 3 summarizes the 9000b78f...md file ---
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = $newGenerateDate

# Widen the "Status" columns now that "Ready for handoff" is longer than
# "In Translation" (Overview!E:F, zh-cn!C, de-de!C).
$overview.Columns.Item(5).ColumnWidth = 16.38
$overview.Columns.Item(6).ColumnWidth = 16.38
$zhcn.Columns.Item(3).ColumnWidth = 16.38
$dede.Columns.Item(3).ColumnWidth = 16.38
